$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 350, pushing existing rows 350-373 down to 352-375
$ws.Rows.Item(350).Resize(2).Insert()

# Row 350 - new data
$ws.Cells.Item(350, 1).Value = 10
$ws.Cells.Item(350, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(350, 3).Value = "La Araucanía"
$ws.Cells.Item(350, 4).Value = 44516
$ws.Cells.Item(350, 5).Value = 9
$ws.Cells.Item(350, 6).Value = "Fruta"
$ws.Cells.Item(350, 7).Value = 100108
$ws.Cells.Item(350, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(350, 9).Value = 100108006
$ws.Cells.Item(350, 10).Value = "Plátano"
$ws.Cells.Item(350, 11).Value = "Barraganete"
$ws.Cells.Item(350, 12).Value = "Primera"
$ws.Cells.Item(350, 13).Value = 35
$ws.Cells.Item(350, 14).Value = 29000
$ws.Cells.Item(350, 15).Value = 29000
$ws.Cells.Item(350, 16).Value = 29000
$ws.Cells.Item(350, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(350, 18).Value = "Ecuador"
$ws.Cells.Item(350, 19).Value = 1450
$ws.Cells.Item(350, 20).Value = 20

# Row 351 - new data
$ws.Cells.Item(351, 1).Value = 10
$ws.Cells.Item(351, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(351, 3).Value = "La Araucanía"
$ws.Cells.Item(351, 4).Value = 44516
$ws.Cells.Item(351, 5).Value = 9
$ws.Cells.Item(351, 6).Value = "Fruta"
$ws.Cells.Item(351, 7).Value = 100108
$ws.Cells.Item(351, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(351, 9).Value = 100108006
$ws.Cells.Item(351, 10).Value = "Plátano"
$ws.Cells.Item(351, 11).Value = "Sin especificar"
$ws.Cells.Item(351, 12).Value = "Pintón"
$ws.Cells.Item(351, 13).Value = 1130
$ws.Cells.Item(351, 14).Value = 17000
$ws.Cells.Item(351, 15).Value = 19000
$ws.Cells.Item(351, 16).Value = 17850
$ws.Cells.Item(351, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(351, 18).Value = "Ecuador"
$ws.Cells.Item(351, 19).Value = 892
$ws.Cells.Item(351, 20).Value = 20
